$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values are written as text (avoid Excel auto-number conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "51.955.74"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "2.861.51"
$ws.Range("E3").Value = "  +2.81%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "351.71"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "113.44"
$ws.Range("E6").Value = "  +4.04%  "
$ws.Range("D7").Value = "0.556"
$ws.Range("E7").Value = "  +1.20%  "
$ws.Range("D8").Value = "0.997"
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("D9").Value = "0.622"
$ws.Range("E9").Value = "  +3.98%  "
$ws.Range("D10").Value = "40.36"
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").Value = "0.0853"
$ws.Range("E12").Value = "  +1.80%  "
$ws.Range("D13").Value = "20.19"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").Value = "7.82"
$ws.Range("E14").Value = "  +2.23%  "
$ws.Range("D15").Value = "3.296.33"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("E16").Value = "  +6.60%  "
$ws.Range("D17").Value = "2.850.02"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").Value = "51.926.77"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("E19").Value = "  +8.65%  "
$ws.Range("D20").Value = "7.66"
$ws.Range("D21").Value = "13.59"
$ws.Range("E21").Value = "  +3.26%  "
$ws.Range("D22").Value = "0.0₃0977"
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("D23").Value = "70.57"
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").Value = "269.90"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("D25").Value = "2.77"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("D26").Value = "26.41"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("D29").Value = "39.18"
$ws.Range("E29").Value = "  +6.58%  "
$ws.Range("D30").Value = "10.59"
$ws.Range("E30").Value = "  +3.47%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "6.31"
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("B32").Value = "OKB"
$ws.Range("C32").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").Value = "52.81"
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").Value = "2.11"
$ws.Range("E33").Value = "  -5.48%  "
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("D35").Value = "0.0894"
$ws.Range("E35").Value = "  +7.57%  "
$ws.Range("D36").Value = "5.62"
$ws.Range("E36").Value = "  +1.32%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "18.91"
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("D39").Value = "3.26"
$ws.Range("E39").Value = "  +3.64%  "
$ws.Range("D40").Value = "2.03"
$ws.Range("E40").Value = "  +2.82%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "2.56"
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.116"
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("D43").Value = "122.63"
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("D44").Value = "22.45"
$ws.Range("E44").Value = "  +1.76%  "
$ws.Range("E45").Value = "  +1.95%  "
$ws.Range("D46").Value = "3.54"
$ws.Range("E46").Value = "  +7.35%  "
$ws.Range("D47").Value = "2.52"
$ws.Range("E47").Value = "  +7.94%  "
$ws.Range("D48").Value = "2.173.22"
$ws.Range("E48").Value = "  +2.08%  "
$ws.Range("D49").Value = "0.245"
$ws.Range("E49").Value = "  +19.39%  "
$ws.Range("D50").Value = "0.963"
$ws.Range("E50").Value = "  +6.13%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").Value = "0.0317"
$ws.Range("E51").Value = "  +12.39%  "

# Restore default style (no explicit numFmt) for column D cells, matching original formatting
$ws.Range("D2:D51").Style = "Normal"

